$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 2, shifting all existing
# data rows (old row 2 -> new row 3, ..., old row 57 -> new row 58).
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting copied from row 1 (the bold
# header row). Reset it back to the plain formatting used by every other
# data row before populating it.
$ws.Range("A2:R2").ClearFormats()

# Column D holds dates; give it the same date number format used by the
# rest of the D column so the new value round-trips as a date like its
# siblings.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with its data.
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = 44515
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 6000
$ws.Range("N2").Value = "`$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 240
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
